$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph. The paragraph needs three runs: an empty leading
#    run, a bold "Meta description" run, and a regular run with the
#    description text - the same run "shape" already used by the duplicated
#    bold-heading paragraph near the end of the document, so we borrow its
#    formatted run structure (via FormattedText) before that paragraph gets
#    removed in a later step.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$templatePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$destRange = $d.Paragraphs.Item(2).Range
$destRange.MoveEnd(1, -1)
$templateRange = $templatePara.Range.Duplicate
$templateRange.MoveEnd(1, -1)
$destRange.FormattedText = $templateRange.FormattedText

# Rename the copied bold run's text to "Meta description".
$renameRange = $d.Paragraphs.Item(2).Range.Duplicate
$renameRange.Find.Execute("Play 6 Wild Sharks Free - Customize Your Wilds & Win Big", $false, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

# Append the (non-bold) description run after "Meta description".
$descRange = $d.Paragraphs.Item(2).Range.Duplicate
$descRange.MoveEnd(1, -1)
$descRange.Collapse(0)
$descRange.InsertAfter(": Read our review of 6 Wild Sharks by 4ThePlayer, featuring customizable wilds, free spins with locked wilds, high-definition graphics, and bonus features. Play now for free.")
$descRange.Bold = 0

# ---------------------------------------------------------------------------
# 2) Remove the duplicated "Play 6 Wild Sharks Free..." (bold) paragraph that
#    used to sit near the end of the document (it was only needed above as a
#    formatting template and the target document no longer contains it).
# ---------------------------------------------------------------------------
$dupHeadingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$dupHeadingPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final paragraph (currently the italic
#    "Read our review ..." meta description) with the new image prompt text,
#    keeping its italic run formatting intact.
# ---------------------------------------------------------------------------
$newPromptText = 'Prompt: Create a cartoon-style feature image for the game "6 Wild Sharks" featuring a happy Maya warrior with glasses. The image should convey the thrill and excitement of hunting sharks while highlighting the game''s unique Wild Choice mechanics and the customizable wild symbols that players can use to increase their chances of winning big. The Maya warrior should be seen holding a fishing rod and standing on a boat with 6 wild sharks jumping out of the water in the background. The image should be bright and colorful, with bold text reading "6 Wild Sharks" and "Customize Your Game" to showcase the game''s innovative features.'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$foundRange = $lastPara.Range.Duplicate
$foundRange.Find.Execute("Read our review of 6 Wild Sharks by 4ThePlayer, featuring customizable wilds, free spins with locked wilds, high-definition graphics, and bonus features. Play now for free.")
$foundRange.Text = $newPromptText

Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Para 1: [" $d.Paragraphs.Item(1).Range.Text "]"
Write-Host "Para 2: [" $d.Paragraphs.Item(2).Range.Text "]"
Write-Host "Para 3: [" $d.Paragraphs.Item(3).Range.Text "]"
Write-Host "Second-to-last Para: [" $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.Text "]"
Write-Host "Last Para: [" $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text "]"
